$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.514.49'
$ws.Range('E2').Value = '  +3.14%  '
$ws.Range('D3').Value = '2.547.63'
$ws.Range('E3').Value = '  +3.27%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '543.39'
$ws.Range('E5').Value = '  +2.11%  '
$ws.Range('D6').Value = '146.76'
$ws.Range('E6').Value = '  +1.82%  '
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').Value = '0.572'
$ws.Range('E8').Value = '  +0.87%  '
$ws.Range('D9').Value = '2.584.21'
$ws.Range('E10').Value = '  +2.93%  '
$ws.Range('E11').Value = '  +1.48%  '
$ws.Range('D12').Value = '5.51'
$ws.Range('E12').Value = '  -1.30%  '
$ws.Range('D13').Value = '0.365'
$ws.Range('E13').Value = '  +4.39%  '
$ws.Range('D14').Value = '2.999.26'
$ws.Range('E14').Value = '  +3.43%  '
$ws.Range('D15').Value = '24.51'
$ws.Range('E15').Value = '  +3.19%  '
$ws.Range('D16').Value = '60.387.63'
$ws.Range('E16').Value = '  +3.16%  '
$ws.Range('E17').Value = '  +5.48%  '
$ws.Range('D18').Value = '2.556.19'
$ws.Range('E18').Value = '  +2.98%  '
$ws.Range('D19').Value = '11.38'
$ws.Range('E19').Value = '  +1.04%  '
$ws.Range('D20').Value = '4.38'
$ws.Range('E20').Value = '  +2.15%  '
$ws.Range('D21').Value = '329.05'
$ws.Range('E21').Value = '  +2.40%  '
$ws.Range('D22').Value = '5.98'
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('D24').Value = '63.04'
$ws.Range('E24').Value = '  +4.12%  '
$ws.Range('D25').Value = '0.443'
$ws.Range('E25').Value = '  +1.56%  '
$ws.Range('E26').Value = '  +4.69%  '
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('D28').Value = '8.09'
$ws.Range('E28').Value = '  +5.43%  '
$ws.Range('D29').Value = '7.20'
$ws.Range('E29').Value = '  +4.45%  '
$ws.Range('D30').Value = '0.0₃0814'
$ws.Range('E30').Value = '  +5.80%  '
$ws.Range('E31').Value = '  +2.77%  '
$ws.Range('E32').Value = '  -0.85%  '
$ws.Range('D33').Value = '164.38'
$ws.Range('E33').Value = '  +4.09%  '
$ws.Range('E34').Value = '  +6.28%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').Value = '18.86'
$ws.Range('E36').Value = '  +2.26%  '
$ws.Range('D37').Value = '4.50'
$ws.Range('E37').Value = '  +3.15%  '
$ws.Range('D38').Value = '1.65'
$ws.Range('E38').Value = '  +3.14%  '
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('D40').Value = '307.73'
$ws.Range('E40').Value = '  +1.10%  '
$ws.Range('D41').Value = '37.19'
$ws.Range('E41').Value = '  +1.93%  '
$ws.Range('D42').Value = '0.846'
$ws.Range('E42').Value = '  +5.15%  '
$ws.Range('D43').Value = '3.78'
$ws.Range('E43').Value = '  +2.25%  '
$ws.Range('D44').Value = '0.613'
$ws.Range('E44').Value = '  +3.80%  '
$ws.Range('D45').Value = '0.990'
$ws.Range('E45').Value = '  -0.49%  '
$ws.Range('D46').Value = '10.85'
$ws.Range('E46').Value = '  +0.72%  '
$ws.Range('D47').Value = '126.97'
$ws.Range('E47').Value = '  +2.43%  '
$ws.Range('D48').Value = '19.16'
$ws.Range('E48').Value = '  +4.19%  '
$ws.Range('D49').Value = '0.0941'
$ws.Range('E49').Value = '  +2.12%  '
$ws.Range('D50').Value = '0.0527'
